$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1099.9048
$ws.Range("J129").Value = 1052.2759
$ws.Range("L129").Value = 3156.8277
$ws.Range("N129").Value = -13156.8277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 69493.16
$ws.Range("J132").Value = 69493.16
$ws.Range("L132").Value = 69493.16
$ws.Range("N132").Value = -79613.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3706500.8
$ws.Range("I31").Value = 2276.8
$ws.Range("J31").Value = 5558612.5
$ws.Range("K31").Value = 2276.8
$ws.Range("L31").Value = 5558612.5
$ws.Range("M31").Value = -1981.8
$ws.Range("N31").Value = -5559202.5
$ws.Range("H34").Value = 3706500.8
$ws.Range("I34").Value = 2276.8
$ws.Range("J34").Value = 5558612.5
$ws.Range("K34").Value = 2276.8
$ws.Range("L34").Value = 5558612.5
$ws.Range("M34").Value = -2074.8
$ws.Range("N34").Value = -5559016.5
$ws.Range("H106").Value = 51666.332
$ws.Range("J106").Value = 51666.332
$ws.Range("L106").Value = 51666.332
$ws.Range("N106").Value = -54190.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 753650.6
$ws.Range("J68").Value = 773582.6
$ws.Range("L68").Value = 2320747.8
$ws.Range("N68").Value = -2322369.8
$ws.Range("H71").Value = 753650.6
$ws.Range("J71").Value = 773582.6
$ws.Range("L71").Value = 6962243.399999999
$ws.Range("N71").Value = -6970355.399999999
$ws.Range("H109").Value = 2403
$ws.Range("I109").Value = 998.1429
$ws.Range("J109").Value = 3222.5
$ws.Range("K109").Value = 2994.4287
$ws.Range("L109").Value = 9667.5
$ws.Range("M109").Value = -1954.4287
$ws.Range("N109").Value = -11747.5
$ws.Range("H114").Value = 2051.875
$ws.Range("I114").Value = 1194
$ws.Range("J114").Value = 3481.6667
$ws.Range("K114").Value = 3582
$ws.Range("L114").Value = 10445.0001
$ws.Range("M114").Value = -328
$ws.Range("N114").Value = -16953.0001
$ws.Range("H117").Value = 3010.6667
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 3010.6667
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 9032.000100000001
$ws.Range("N117").Value = -15916.0001
$ws.Range("M117").ClearContents()
$ws.Range("H125").Value = 2311101
$ws.Range("I125").Value = 7501607.5
$ws.Range("J125").Value = 4209.222
$ws.Range("K125").Value = 22504822.5
$ws.Range("L125").Value = 12627.666
$ws.Range("M125").Value = -22499902.5
$ws.Range("N125").Value = -22467.666
$ws.Range("H129").Value = 116959
$ws.Range("I129").Value = 334154
$ws.Range("J129").Value = 1973.4117
$ws.Range("K129").Value = 1002462
$ws.Range("L129").Value = 5920.2351
$ws.Range("M129").Value = -997462
$ws.Range("N129").Value = -15920.2351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 45607.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 45607.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 45607.5
$ws.Range("N125").Value = -55447.5
$ws.Range("H127").Value = 44245
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 44245
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 44245
$ws.Range("N127").Value = -54165
$ws.Range("H128").Value = 28429
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 28429
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 28429
$ws.Range("N128").Value = -38389
$ws.Range("H129").Value = 45429
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45429
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45429
$ws.Range("N129").Value = -55429
$ws.Range("H130").Value = 48296
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48296
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48296
$ws.Range("N130").Value = -58336
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 3501.4827
$ws.Range("I132").Value = 2564.742
$ws.Range("J132").Value = 4577
$ws.Range("K132").Value = 7694.226000000001
$ws.Range("L132").Value = 13731
$ws.Range("M132").Value = -5164.226000000001
$ws.Range("N132").Value = -18791
$ws.Range("H133").Value = 43313.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 43313.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 43313.25
$ws.Range("N133").Value = -48373.25
$ws.Range("H134").Value = 69199.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 69199.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 69199.5
$ws.Range("N134").Value = -79339.5
$ws.Range("H135").Value = 37843.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 37843.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 37843.332
$ws.Range("N135").Value = -47983.332
$ws.Range("H136").Value = 2898.7058
$ws.Range("I136").Value = 2058.1667
$ws.Range("J136").Value = 4916
$ws.Range("K136").Value = 6174.500100000001
$ws.Range("L136").Value = 14748
$ws.Range("M136").Value = -3624.500100000001
$ws.Range("N136").Value = -19848
$ws.Range("H137").Value = 44662.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 44662.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 44662.5
$ws.Range("N137").Value = -54862.5
$ws.Range("H138").Value = 58000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 58000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 58000
$ws.Range("N138").Value = -68280
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 39999.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39999.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39999.332
$ws.Range("N140").Value = -50359.332
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 40869.7
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 40869.7
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 40869.7
$ws.Range("N119").Value = -50545.7
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 30420
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 30420
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 30420
$ws.Range("N121").Value = -33914
$ws.Range("H122").Value = 1787739.5
$ws.Range("I122").Value = 3176936.2
$ws.Range("J122").Value = 1629.2858
$ws.Range("K122").Value = 9530808.600000001
$ws.Range("L122").Value = 4887.857400000001
$ws.Range("M122").Value = -9528358.600000001
$ws.Range("N122").Value = -9787.8574
$ws.Range("H123").Value = 43076.332
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 43076.332
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 43076.332
$ws.Range("N123").Value = -52876.332
$ws.Range("H124").Value = 45429
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 45429
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 45429
$ws.Range("N124").Value = -55249
$ws.Range("H125").Value = 36633.332
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 36633.332
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 36633.332
$ws.Range("N125").Value = -46473.332
$ws.Range("H126").Value = 1402243.1
$ws.Range("I126").Value = 1635583.6
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 4906750.800000001
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -4904280.800000001
$ws.Range("N126").Value = -11540
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 48000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 48000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960
$ws.Range("H129").Value = 40429
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 40429
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 40429
$ws.Range("N129").Value = -50429
$ws.Range("H130").Value = 30428.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 30428.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 30428.5
$ws.Range("N130").Value = -40468.5
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 1209782.9
$ws.Range("I132").Value = 1673740.1
$ws.Range("J132").Value = 3494.1
$ws.Range("K132").Value = 5021220.300000001
$ws.Range("L132").Value = 10482.3
$ws.Range("M132").Value = -5018690.300000001
$ws.Range("N132").Value = -15542.3
$ws.Range("H133").Value = 50531.168
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50531.168
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50531.168
$ws.Range("N133").Value = -60651.168
$ws.Range("H135").Value = 38990.266
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38990.266
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38990.266
$ws.Range("N135").Value = -49130.266
$ws.Range("H136").Value = 584440.94
$ws.Range("I136").Value = 898320.75
$ws.Range("J136").Value = 1521.2142
$ws.Range("K136").Value = 2694962.25
$ws.Range("L136").Value = 4563.642599999999
$ws.Range("M136").Value = -2692412.25
$ws.Range("N136").Value = -9663.6426
$ws.Range("H137").Value = 47000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 47000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 47000
$ws.Range("N137").Value = -57200
$ws.Range("H138").Value = 40997.145
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40997.145
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40997.145
$ws.Range("N138").Value = -51277.145
$ws.Range("H139").Value = 47960
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47960
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47960
$ws.Range("N139").Value = -58240
$ws.Range("H140").Value = 41586.555
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 41586.555
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 41586.555
$ws.Range("N140").Value = -51946.555
$ws.Range("H141").Value = 39200
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39200
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 39200
$ws.Range("N141").Value = -49560
